# Update "想去人数" (want-to-go count) figures for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.
#   - 南宁·布谷鸟动漫展5th              (row 2): 490  -> 492
#   - 南宁·2024良牙动漫秋季盛典（秋典） (row 3): 5812 -> 5831

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 492
    $ws.Range("F3").Value = 5831
}
